# Update version string from "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
# to "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$wb = $excel.ActiveWorkbook

$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

# A2: "Version: ..."
$wsAbout.Range("A2").Value = "Version: " + $newVersion

# A6: Recommended Citation containing the version string
$newCitation = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for Bulga Coal Mine, Australia, M0017, version ''' + $newVersion + '''. (See the CC license for attribution requirements if sharing or adapting the data set.)'
$wsAbout.Range("A6").Value = $newCitation

# S2:S37 on "Boundaries and methane sources" sheet hold the build_version value
for ($r = 2; $r -le 37; $r++) {
    $wsData.Cells.Item($r, 19).Value = $newVersion
}
